# Applies the cryptos-list refresh described by the diff (commit: 'Updated cryptos list ... with GitHub Actions').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.059.55'
$ws.Range('D3').Value = '2.749.76'
$ws.Range('E3').Value = '  +3.76%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '606.24'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '167.19'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.68%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').Value = '2.748.77'
$ws.Range('E9').Value = '  +3.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.142'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.365'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.79%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.36'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.93'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.00%  '
$ws.Range('D15').Value = '3.245.59'
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '68.968.59'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '2.729.14'
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.98'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.83%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.74'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.53%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '369.68'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('E22').Value = '  +3.20%  '
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.18'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('E28').Value = '  +3.34%  '
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '599.60'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +7.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('E35').Value = '  +3.34%  '
$ws.Range('E36').Value = '  +4.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = "Normal"
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '163.40'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '20.17'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.52'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.91%  '
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '18.02'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('D45').Value = '0.0₆0319'
$ws.Range('E45').Value = '  -5.04%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '159.11'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('E48').Value = '  +5.41%  '
$ws.Range('E49').Value = '  +6.81%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.611'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.88%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '22.15'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.54%  '
